# Creado grafico de tipos de modelo
# Insert a new "MAE" column before the existing "Tipo" column.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column at D, shifting the existing "Tipo"/"single" column to E.
$ws.Columns.Item(4).Insert()

# Populate the new column's header and value.
$ws.Range("D1").Value = "MAE"
$ws.Range("D2").Value = 1.29915171380136
